$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1192057.70884408
$ws.Range("C2").Value = 167978496243.219
$ws.Range("D2").Value = 8087276.1799207

$ws.Range("B3").Value = 1619354.83040346
$ws.Range("C3").Value = 1675532.91028309
$ws.Range("D3").Value = 1622785.81590323
